$wb = $excel.ActiveWorkbook

$wsResults = $wb.Worksheets.Item("results")
$wsStats = $wb.Worksheets.Item("stats")

$wsResults.Range("B2").Value = 59.452
$wsResults.Range("C2").Value = 50.59199999999999
$wsResults.Range("D2").Value = 50.59199999999999
$wsResults.Range("E2").Value = 50.59199999999999
$wsResults.Range("F2").Value = 50.59199999999999
$wsResults.Range("G2").Value = 0.1192222296979077
$wsResults.Range("B3").Value = 20.726
$wsResults.Range("C3").Value = 42.038
$wsResults.Range("D3").Value = 42.038
$wsResults.Range("E3").Value = 42.038
$wsResults.Range("F3").Value = 42.038
$wsResults.Range("G3").Value = 0.8226189327414839
$wsResults.Range("B4").Value = 61.796
$wsResults.Range("C4").Value = 54.55399999999999
$wsResults.Range("D4").Value = 54.554
$wsResults.Range("E4").Value = 54.554
$wsResults.Range("F4").Value = 54.554
$wsResults.Range("G4").Value = 0.09375364101236328
$wsResults.Range("B5").Value = 40.108
$wsResults.Range("C5").Value = 37.14
$wsResults.Range("D5").Value = 37.14
$wsResults.Range("E5").Value = 37.14
$wsResults.Range("F5").Value = 37.14
$wsResults.Range("G5").Value = 0.05920015956916314
$wsResults.Range("B6").Value = 41.554
$wsResults.Range("C6").Value = 65.79599999999999
$wsResults.Range("D6").Value = 65.79599999999999
$wsResults.Range("E6").Value = 65.79599999999999
$wsResults.Range("F6").Value = 65.79599999999999
$wsResults.Range("G6").Value = 0.4667083794580542
$wsStats.Range("C2").Value = 1840
$wsStats.Range("D2").Value = 0.02247039030771703
$wsStats.Range("E2").Value = 0.7443019769852981
$wsStats.Range("F2").Value = 1840
$wsStats.Range("G2").Value = 0.05537386040668935
$wsStats.Range("H2").Value = 0.4904221027391031
$wsStats.Range("I2").Value = 0.03413118747994304
$wsStats.Range("J2").Value = 0.08356877986807376
$wsStats.Range("K2").Value = 0.01999597204849124
$wsStats.Range("C3").Value = 866
$wsStats.Range("D3").Value = 0.0009043698664754629
$wsStats.Range("E3").Value = 0.3264306769706309
$wsStats.Range("F3").Value = 866
$wsStats.Range("G3").Value = 0.02500208560377359
$wsStats.Range("H3").Value = 0.2106004991801456
$wsStats.Range("I3").Value = 0.007719359477050602
$wsStats.Range("J3").Value = 0.04627484793309122
$wsStats.Range("K3").Value = 0.009199967258609831
$wsStats.Range("C4").Value = 866
$wsStats.Range("D4").Value = 0.009958986309356987
$wsStats.Range("E4").Value = 0.3579659190727398
$wsStats.Range("F4").Value = 866
$wsStats.Range("G4").Value = 0.02703198709059507
$wsStats.Range("H4").Value = 0.2321864258265123
$wsStats.Range("I4").Value = 0.009908543550409377
$wsStats.Range("J4").Value = 0.04946563148405403
$wsStats.Range("K4").Value = 0.009843050269410014
$wsStats.Range("C5").Value = 866
$wsStats.Range("D5").Value = 0.001394452061504126
$wsStats.Range("E5").Value = 0.3342486170586199
$wsStats.Range("F5").Value = 866
$wsStats.Range("G5").Value = 0.02607315022032708
$wsStats.Range("H5").Value = 0.2134779418120161
$wsStats.Range("I5").Value = 0.008885901304893196
$wsStats.Range("J5").Value = 0.04743586608674377
$wsStats.Range("K5").Value = 0.009680040297098458
$wsStats.Range("C6").Value = 4805
$wsStats.Range("E6").Value = 0.9688141721999273
$wsStats.Range("C7").Value = 1563
$wsStats.Range("D7").Value = 0.01898405444808304
$wsStats.Range("E7").Value = 0.6381196749862283
$wsStats.Range("F7").Value = 1563
$wsStats.Range("G7").Value = 0.04755687050055712
$wsStats.Range("H7").Value = 0.4147311503766105
$wsStats.Range("I7").Value = 0.03226113377604634
$wsStats.Range("J7").Value = 0.07299588446039706
$wsStats.Range("K7").Value = 0.01736028608866036
$wsStats.Range("C8").Value = 707
$wsStats.Range("D8").Value = 0.000774587388150394
$wsStats.Range("E8").Value = 0.2742936819558963
$wsStats.Range("F8").Value = 707
$wsStats.Range("G8").Value = 0.0207951336633414
$wsStats.Range("H8").Value = 0.1727409242885187
$wsStats.Range("I8").Value = 0.008141346508637071
$wsStats.Range("J8").Value = 0.04147404141258448
$wsStats.Range("K8").Value = 0.007689272402785718
$wsStats.Range("C9").Value = 707
$wsStats.Range("D9").Value = 0.009927223902195692
$wsStats.Range("E9").Value = 0.3200639389688149
$wsStats.Range("F9").Value = 707
$wsStats.Range("G9").Value = 0.02328084665350616
$wsStats.Range("H9").Value = 0.2040852522477508
$wsStats.Range("I9").Value = 0.01169244933407754
$wsStats.Range("J9").Value = 0.04653520707506686
$wsStats.Range("K9").Value = 0.008607124909758568
$wsStats.Range("C10").Value = 707
$wsStats.Range("D10").Value = 0.001186518464237452
$wsStats.Range("E10").Value = 0.2764706989983097
$wsStats.Range("F10").Value = 707
$wsStats.Range("G10").Value = 0.02072662871796638
$wsStats.Range("H10").Value = 0.1732716474216431
$wsStats.Range("I10").Value = 0.009598656557500362
$wsStats.Range("J10").Value = 0.04189795663114637
$wsStats.Range("K10").Value = 0.007737652049399912
$wsStats.Range("C11").Value = 4805
$wsStats.Range("E11").Value = 1.029153865994886
$wsStats.Range("C12").Value = 1823
$wsStats.Range("D12").Value = 0.02205983060412109
$wsStats.Range("E12").Value = 0.726483765989542
$wsStats.Range("F12").Value = 1823
$wsStats.Range("G12").Value = 0.0547923871781677
$wsStats.Range("H12").Value = 0.4769193021347746
$wsStats.Range("I12").Value = 0.03116031002718955
$wsStats.Range("J12").Value = 0.08260764647275209
$wsStats.Range("K12").Value = 0.02003651857376099
$wsStats.Range("C13").Value = 855
$wsStats.Range("D13").Value = 0.0008611893281340599
$wsStats.Range("E13").Value = 0.3242480909684673
$wsStats.Range("F13").Value = 855
$wsStats.Range("G13").Value = 0.02528036874718964
$wsStats.Range("H13").Value = 0.2062780619598925
$wsStats.Range("I13").Value = 0.007863161503337324
$wsStats.Range("J13").Value = 0.04654615197796375
$wsStats.Range("K13").Value = 0.009509829105809331
$wsStats.Range("C14").Value = 855
$wsStats.Range("D14").Value = 0.00881702231708914
$wsStats.Range("E14").Value = 0.3386389480438083
$wsStats.Range("F14").Value = 855
$wsStats.Range("G14").Value = 0.02580453595146537
$wsStats.Range("H14").Value = 0.217749941861257
$wsStats.Range("I14").Value = 0.009843258303590119
$wsStats.Range("J14").Value = 0.04759739525616169
$wsStats.Range("K14").Value = 0.009342187666334212
$wsStats.Range("C15").Value = 855
$wsStats.Range("D15").Value = 0.001334894681349397
$wsStats.Range("E15").Value = 0.3278302811086178
$wsStats.Range("F15").Value = 855
$wsStats.Range("G15").Value = 0.02561569621320814
$wsStats.Range("H15").Value = 0.2084406393114477
$wsStats.Range("I15").Value = 0.008960254141129553
$wsStats.Range("J15").Value = 0.04726952570490539
$wsStats.Range("K15").Value = 0.009454978513531387
$wsStats.Range("C16").Value = 4805
$wsStats.Range("E16").Value = 0.9926358649972826
$wsStats.Range("C17").Value = 1245
$wsStats.Range("D17").Value = 0.01364321040455252
$wsStats.Range("E17").Value = 0.4918821790488437
$wsStats.Range("F17").Value = 1245
$wsStats.Range("G17").Value = 0.03681550000328571
$wsStats.Range("H17").Value = 0.3218109229346737
$wsStats.Range("I17").Value = 0.02245280146598816
$wsStats.Range("J17").Value = 0.0559447273844853
$wsStats.Range("K17").Value = 0.01334235642571002
$wsStats.Range("C18").Value = 591
$wsStats.Range("D18").Value = 0.0006267315475270152
$wsStats.Range("E18").Value = 0.2261388130718842
$wsStats.Range("F18").Value = 591
$wsStats.Range("G18").Value = 0.01741581899113953
$wsStats.Range("H18").Value = 0.1438372770790011
$wsStats.Range("I18").Value = 0.005481239408254623
$wsStats.Range("J18").Value = 0.0332215492380783
$wsStats.Range("K18").Value = 0.0064511300297454
$wsStats.Range("C19").Value = 591
$wsStats.Range("D19").Value = 0.006310588214546442
$wsStats.Range("E19").Value = 0.2399847289780155
$wsStats.Range("F19").Value = 591
$wsStats.Range("G19").Value = 0.01770824962295592
$wsStats.Range("H19").Value = 0.1552030059974641
$wsStats.Range("I19").Value = 0.007016404648311436
$wsStats.Range("J19").Value = 0.03406144538894296
$wsStats.Range("K19").Value = 0.006488861632533371
$wsStats.Range("C20").Value = 591
$wsStats.Range("D20").Value = 0.0009861706057563424
$wsStats.Range("E20").Value = 0.2281700660241768
$wsStats.Range("F20").Value = 591
$wsStats.Range("G20").Value = 0.01733138435520232
$wsStats.Range("H20").Value = 0.1450486677931622
$wsStats.Range("I20").Value = 0.006312734563834965
$wsStats.Range("J20").Value = 0.03404790908098221
$wsStats.Range("K20").Value = 0.006336680613458157
$wsStats.Range("C21").Value = 4805
$wsStats.Range("E21").Value = 1.202622185926884
$wsStats.Range("C22").Value = 2170
$wsStats.Range("D22").Value = 0.02668864489533007
$wsStats.Range("E22").Value = 0.8560476240236312
$wsStats.Range("F22").Value = 2170
$wsStats.Range("G22").Value = 0.06420882721431553
$wsStats.Range("H22").Value = 0.5605792781570926
$wsStats.Range("I22").Value = 0.04085709701757878
$wsStats.Range("J22").Value = 0.09624304971657693
$wsStats.Range("K22").Value = 0.02330316929146647
$wsStats.Range("C23").Value = 976
$wsStats.Range("D23").Value = 0.0008980858838185668
$wsStats.Range("E23").Value = 0.3461714800214395
$wsStats.Range("F23").Value = 976
$wsStats.Range("G23").Value = 0.02731961489189416
$wsStats.Range("H23").Value = 0.2216545251430944
$wsStats.Range("I23").Value = 0.009241797495633364
$wsStats.Range("J23").Value = 0.04734734527301043
$wsStats.Range("K23").Value = 0.01030507707037032
$wsStats.Range("C24").Value = 976
$wsStats.Range("D24").Value = 0.01032575033605099
$wsStats.Range("E24").Value = 0.3793565769447014
$wsStats.Range("F24").Value = 976
$wsStats.Range("G24").Value = 0.02844432753045112
$wsStats.Range("H24").Value = 0.2465887044090778
$wsStats.Range("I24").Value = 0.01109551044646651
$wsStats.Range("J24").Value = 0.05100856814533472
$wsStats.Range("K24").Value = 0.01063325954601169
$wsStats.Range("C25").Value = 976
$wsStats.Range("D25").Value = 0.001531971851363778
$wsStats.Range("E25").Value = 0.3688036509556696
$wsStats.Range("F25").Value = 976
$wsStats.Range("G25").Value = 0.02935354050714523
$wsStats.Range("H25").Value = 0.2364394528558478
$wsStats.Range("I25").Value = 0.01073132548481226
$wsStats.Range("J25").Value = 0.05059266509488225
$wsStats.Range("K25").Value = 0.01033912471029907
$wsStats.Range("C26").Value = 4805
$wsStats.Range("E26").Value = 0.841553280246444
$wsStats.Range("C27").Value = 1728.2
$wsStats.Range("D27").Value = 0.02076922613196075
$wsStats.Range("E27").Value = 0.6913670442067087
$wsStats.Range("F27").Value = 1728.2
$wsStats.Range("G27").Value = 0.05174948906060308
$wsStats.Range("H27").Value = 0.4528925512684509
$wsStats.Range("I27").Value = 0.03217250595334917
$wsStats.Range("J27").Value = 0.07827201758045703
$wsStats.Range("K27").Value = 0.01880766048561782
$wsStats.Range("C28").Value = 799
$wsStats.Range("D28").Value = 0.0008129928028210998
$wsStats.Range("E28").Value = 0.2994565485976636
$wsStats.Range("F28").Value = 799
$wsStats.Range("G28").Value = 0.02316260437946767
$wsStats.Range("H28").Value = 0.1910222575301304
$wsStats.Range("I28").Value = 0.007689380878582596
$wsStats.Range("J28").Value = 0.04297278716694564
$wsStats.Range("K28").Value = 0.00863105517346412
$wsStats.Range("C29").Value = 799
$wsStats.Range("D29").Value = 0.00906791421584785
$wsStats.Range("E29").Value = 0.327202022401616
$wsStats.Range("F29").Value = 799
$wsStats.Range("G29").Value = 0.02445398936979472
$wsStats.Range("H29").Value = 0.2111626660684124
$wsStats.Range("I29").Value = 0.009911233256570996
$wsStats.Range("J29").Value = 0.04573364946991205
$wsStats.Range("K29").Value = 0.00898289680480957
$wsStats.Range("C30").Value = 799
$wsStats.Range("D30").Value = 0.001286801532842219
$wsStats.Range("E30").Value = 0.3071046628290787
$wsStats.Range("F30").Value = 799
$wsStats.Range("G30").Value = 0.02382008000276983
$wsStats.Range("H30").Value = 0.1953356698388234
$wsStats.Range("I30").Value = 0.008897774410434068
$wsStats.Range("J30").Value = 0.044248784519732
$wsStats.Range("K30").Value = 0.008709695236757398
$wsStats.Range("C31").Value = 4805
$wsStats.Range("E31").Value = 1.006955873873085
